# Apply the "APNLP_templates" update:
#  - Add 10 new worksheets (one per ConceptNet-style relation/edge type),
#    each populated with a small table of "Template"/"Expected POS" rows,
#    so that the app can randomly choose among several templates per relation.
#  - Keep the existing "Tabelle1" overview sheet as-is (only the saved
#    cell-selection on it changes).

$wb = $excel.ActiveWorkbook
$tabelle1 = $wb.Worksheets.Item(1)

function Add-TemplateSheet {
    param(
        [string]$Name,
        [object]$After,
        [string[][]]$Rows
    )

    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $After)
    $ws.Name = $Name

    for ($r = 0; $r -lt $Rows.Length; $r++) {
        $row = $Rows[$r]
        for ($c = 0; $c -lt $row.Length; $c++) {
            $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
        }
    }

    # Bold header row (row 1), matching the existing "Tabelle1" style.
    $headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $Rows[0].Length))
    $headerRange.Font.Bold = $true

    return $ws
}

$prev = $tabelle1

$hasARows = @(
    , @("Template", "Expected POS")
    , @("It has <HasA>", "NOUN")
)
$prev = Add-TemplateSheet "HasA" $prev $hasARows

$hasPropertyRows = @(
    , @("Template", "Expected POS")
    , @("It is <HasProperty>", "NOUN, VERB, ADJ")
)
$prev = Add-TemplateSheet "HasProperty" $prev $hasPropertyRows

$partOfRows = @(
    , @("Template", "Expected POS")
    , @("It is part of <PartOf>", "NOUN")
)
$prev = Add-TemplateSheet "PartOf" $prev $partOfRows

$madeOfRows = @(
    , @("Template", "Expected POS")
    , @("It is made of <MadeOf>", "NOUN")
)
$prev = Add-TemplateSheet "MadeOf" $prev $madeOfRows

$atLocationRows = @(
    , @("Template", "Expected POS")
    , @("It can usually be found in <AtLocation>", "NOUN")
    , @("Its favorite spot is <AtLocation>", "NOUN")
    , @("It likes to be at <AtLocation>", "NOUN")
    , @("In summer, it likes to go to <AtLocation>", "NOUN")
    , @("It likes to hide in <AtLocation>", "NOUN")
)
$prev = Add-TemplateSheet "AtLocation" $prev $atLocationRows

$capableOfRows = @(
    , @("Template", "Expected POS")
    , @("Its special ability is to <CapableOf>", "VERB")
    , @("It can <CapableOf>", "VERB")
    , @("It is capable of <CapableOf>", "VERB")
    , @("Its favorite thing to do is <CapableOf>", "VERB")
    , @("Its specialty lies in <CapableOf>", "VERB")
)
$prev = Add-TemplateSheet "CapableOf" $prev $capableOfRows

$notDesiresRows = @(
    , @("Template", "Expected POS")
    , @("It hates <NotDesires>", "NOUN, VERB")
    , @("It does not desire <NotDesires>", "NOUN, VERB")
    , @("It strongly dislikes <NotDesires>", "NOUN, VERB")
    , @("It is afraid of <NotDesires>", "NOUN, VERB")
    , @("Its least favorite thing is <NotDesires>", "NOUN, VERB")
)
$prev = Add-TemplateSheet "NotDesires" $prev $notDesiresRows

$desiresRows = @(
    , @("Template", "Expected POS")
    , @("It loves <Desires>", "NOUN, VERB")
    , @("It likes <Desires>", "NOUN, VERB")
    , @("It desires <Desires>", "NOUN, VERB")
    , @("Its favorite thing is <Desires>", "NOUN, VERB")
    , @("It lives for <Desires>", "NOUN, VERB")
)
$prev = Add-TemplateSheet "Desires" $prev $desiresRows

$usedForRows = @(
    , @("Template", "Expected POS")
    , @("It is used for <UsedFor>", "VERB")
    , @("It is especially good at <UsedFor>", "VERB")
    , @("Its strength lies in <UsedFor>", "VERB")
)
$prev = Add-TemplateSheet "UsedFor" $prev $usedForRows

$isARows = @(
    , @("Template", "Expected POS")
    , @("It is a <IsA>", "NOUN")
)
$prev = Add-TemplateSheet "IsA" $prev $isARows

# Restore focus to the original overview sheet and its (updated) selection.
$tabelle1.Activate() | Out-Null
$tabelle1.Range("B12").Select() | Out-Null
